$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New results rows for the GoodSplit algorithm run against TeacherRL
# (row 33: plain maxLen:2 CE processing, row 34: maxLen:2 + EQtoStop)

$ws.Range("A33").Value = 1
$ws.Range("B33").Value = 6315
$ws.Range("C33").Value = 6315
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 46043
$ws.Range("F33").Value = "../data/tests/sequences/Moore_R10_PDS.fsm"
$ws.Range("G33").Value = "GoodSplit"
$ws.Range("H33").Value = "maxLen:2"
$ws.Range("I33").Value = "TeacherRL"

$ws.Range("A34").Value = 1
$ws.Range("B34").Value = 3069
$ws.Range("C34").Value = 3069
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 21309
$ws.Range("F34").Value = "../data/tests/sequences/Moore_R10_PDS.fsm"
$ws.Range("G34").Value = "GoodSplit"
$ws.Range("H34").Value = "maxLen:2 + EQtoStop"
$ws.Range("I34").Value = "TeacherRL"

$ws.Range("E34").Select()
